$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("相談件数")

# The sheet currently ends with a data row (row 109) followed by a note
# row (row 110, "※4/8より..."). A new day of data needs to be inserted
# right before that note row, so: insert a fresh row at 110 (this pushes
# the existing note row down to 111), then populate the new row 110 with
# the next day's figures.
$ws.Rows(110).Insert()

$ws.Range("A110").Value = 43965
$ws.Range("B110").Value = 229
$ws.Range("C110").Value = 37074
$ws.Range("D110").Value = 42
$ws.Range("E110").Value = 7530

# Keep the active selection in sync with the new last data row, matching
# the saved selection state in the sheet.
$ws.Range("A111").Select()

# Extend the print area by one row to include the newly added data row.
$ws.PageSetup.PrintArea = '$A$1:$E$112'
